$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = 'Cluster Name'
$ws.Range("B1").Value = 'Actives Cases'

# Update data rows (cluster name + active cases)
$ws.Range("A2").Value = '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North'
$ws.Range("B2").Value = 14
$ws.Range("A3").Value = '3398 BlueCross Elly Kay Mordialloc'
$ws.Range("B3").Value = 41
$ws.Range("A4").Value = '3601 Baptcare Westhaven community'
$ws.Range("B4").Value = 21
$ws.Range("A5").Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Range("B5").Value = 23
$ws.Range("A6").Value = '3939 Bupa Aged Care Eastwood'
$ws.Range("B6").Value = 11
$ws.Range("A7").Value = '3975 Aurrum Aged Care Brunswick West'
$ws.Range("B7").Value = 11
$ws.Range("A8").Value = '3988 Kerala Manor Aged Care Diamond Creek'
$ws.Range("B8").Value = 10
$ws.Range("A9").Value = '4257 BlueCross The Gables Camberwell'
$ws.Range("B9").Value = 28
$ws.Range("A10").Value = '4295 Hope Aged Care Sunshine West'
$ws.Range("B10").Value = 31
$ws.Range("A11").Value = '44087 Fitzroy Primary School Fitzroy'
$ws.Range("B11").Value = 22
$ws.Range("A12").Value = '44098 Stawell Primary School'
$ws.Range("B12").Value = 27
$ws.Range("A13").Value = '44366 Lysterfield Primary School Lysterfield'
$ws.Range("B13").Value = 12
$ws.Range("A14").Value = '44444 Nar Nar Goon Primary School Nar Nar Goon'
$ws.Range("B14").Value = 18
$ws.Range("A15").Value = '44630 Black Rock Primary School Black Rock'
$ws.Range("B15").Value = 22
$ws.Range("A16").Value = '44666 Gardenvale Primary School Brighton East'
$ws.Range("B16").Value = 21
$ws.Range("A17").Value = '44811 Dandenong North Primary School Dandenong'
$ws.Range("B17").Value = 20
$ws.Range("A18").Value = '44865 Parktone Primary School Parkdale'
$ws.Range("B18").Value = 22
$ws.Range("A19").Value = '44950 Templestowe Valley Primary School Templestowe Lower'
$ws.Range("B19").Value = 65
$ws.Range("A20").Value = '44982 Diamond Creek East Primary School Diamond Creek'
$ws.Range("B20").Value = 11
$ws.Range("A21").Value = '45026 Churchill North Primary School Churchill'
$ws.Range("B21").Value = 15
$ws.Range("A22").Value = '45248 Brookside P-9 College Caroline Springs'
$ws.Range("B22").Value = 31
$ws.Range("A23").Value = '45249 Creekside K-9 College Caroline Springs'
$ws.Range("B23").Value = 12
$ws.Range("A24").Value = '45267 Epping Views Primary School Epping'
$ws.Range("B24").Value = 20
$ws.Range("A25").Value = '45315 Red Hill Consolidated School Red Hill'
$ws.Range("B25").Value = 15
$ws.Range("A26").Value = '45518 Ashwood High School Ashwood'
$ws.Range("B26").Value = 21
$ws.Range("A27").Value = '45569 Nhill College Nhill'
$ws.Range("B27").Value = 33
$ws.Range("A28").Value = '45585 Mount Ridley College Craigieburn'
$ws.Range("B28").Value = 11
$ws.Range("A29").Value = '45648 St Brendans Primary School Shepparton'
$ws.Range("B29").Value = 33
$ws.Range("A30").Value = '4574 Village Glen Aged Care Residences Mornington'
$ws.Range("B30").Value = 11
$ws.Range("A31").Value = '45755 St Patricks Catholic Parish Primary School Mentone'
$ws.Range("B31").Value = 17
$ws.Range("A32").Value = '45784 Holy Rosary Primary School White Hills'
$ws.Range("B32").Value = 26
$ws.Range("A33").Value = '45810 Loreto Mandeville Hall Toorak'
$ws.Range("B33").Value = 11
$ws.Range("A34").Value = '45846 St Mary''s School Mooroopna'
$ws.Range("B34").Value = 18
$ws.Range("A35").Value = '45848 St Kevin''s College Toorak Glendalough Campus Junior School'
$ws.Range("B35").Value = 16
$ws.Range("A36").Value = '45912 St Bernadette''s Catholic Primary School Sunshine North Exposure Site'
$ws.Range("B36").Value = 10
$ws.Range("A37").Value = '45950 St Luke''s Primary School Lalor'
$ws.Range("B37").Value = 21
$ws.Range("A38").Value = '46028 St Anne''s Catholic Primary School Sunbury'
$ws.Range("B38").Value = 12
$ws.Range("A39").Value = '46037 Nazareth Catholic Primary School Grovedale'
$ws.Range("B39").Value = 27
$ws.Range("A40").Value = '46050 Our Lady''s Catholic Primary School Craigieburn'
$ws.Range("B40").Value = 12
$ws.Range("A41").Value = '46052 St. Francis of Assisi Primary School Mill Park'
$ws.Range("B41").Value = 26
$ws.Range("A42").Value = '46093 St Brendan''s Primary School Somerville'
$ws.Range("B42").Value = 14
$ws.Range("A43").Value = '46095 Bethany Catholic Primary School Werribee'
$ws.Range("B43").Value = 11
$ws.Range("A44").Value = '46105 Christ the Priest Primary School Caroline Springs'
$ws.Range("B44").Value = 41
$ws.Range("A45").Value = '46115 St Luke''s Catholic Primary School Shepparton North'
$ws.Range("B45").Value = 10
$ws.Range("A46").Value = '46117 Marymede Catholic College South Morang'
$ws.Range("B46").Value = 13
$ws.Range("A47").Value = '46125 Our Lady of the Southern Cross Primary School Manor Lakes'
$ws.Range("B47").Value = 37
$ws.Range("A48").Value = '46221 Bialik College Hawthorn'
$ws.Range("B48").Value = 13
$ws.Range("A49").Value = '46239 Gilson College Taylors Hill'
$ws.Range("B49").Value = 12
$ws.Range("A50").Value = '46287 Oakleigh Grammar Melbourne Private School Oakleigh'
$ws.Range("B50").Value = 25
$ws.Range("A51").Value = '46328 Ilim College Inverloch Crescent Dallas'
$ws.Range("B51").Value = 10
$ws.Range("A52").Value = '46390 Al Siraat College Epping'
$ws.Range("B52").Value = 30
$ws.Range("A53").Value = '50584 St Mary of the Cross MacKillop Primary School Epping'
$ws.Range("B53").Value = 13
$ws.Range("A54").Value = '51529 Sirius College Primary School Dallas'
$ws.Range("B54").Value = 14
$ws.Range("A55").Value = 'Alfred Health Caulfield Hospital'
$ws.Range("B55").Value = 10
$ws.Range("A56").Value = 'Alfred Health The Alfred Hospital Melbourne'
$ws.Range("B56").Value = 13
$ws.Range("A57").Value = 'Covenant College Bell Post Hill'
$ws.Range("B57").Value = 25
$ws.Range("A58").Value = 'Epping Views Primary School Camp Cape Schanck'
$ws.Range("B58").Value = 14
$ws.Range("A59").Value = 'Hamilton Country Music Festival Hamilton Golf Club Hamilton'
$ws.Range("B59").Value = 12
$ws.Range("A60").Value = 'House Party 27 November Private Residence Brunswick West'
$ws.Range("B60").Value = 25
$ws.Range("A61").Value = 'Islamic College of Melbourne Tarneit'
$ws.Range("B61").Value = 16
$ws.Range("A62").Value = 'Social Gathering 20 November Sunbury'
$ws.Range("B62").Value = 20
$ws.Range("A63").Value = 'Springside Primary School Caroline Springs'
$ws.Range("B63").Value = 25
$ws.Range("A64").Value = 'St Josephs Catholic Primary School Warragul'
$ws.Range("B64").Value = 13
$ws.Range("A65").Value = 'St Vincents Hospital Melbourne Emergency Department Fitzroy'
$ws.Range("B65").Value = 12
$ws.Range("A66").Value = 'The Village Early Learning Centre Sandringham'
$ws.Range("B66").Value = 15
$ws.Range("A67").Value = 'Wagstaff Meat Processing Plant Cranbourne East'
$ws.Range("B67").Value = 35
$ws.Range("A68").Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Range("B68").Value = 15
$ws.Range("A69").Value = 'Western Health Sunshine Hospital Emergency Department St Albans'
$ws.Range("B69").Value = 11
